$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 4462.9536
$ws.Range("I15").Value = 4462.9536
$ws.Range("K15").Value = 13388.8608
$ws.Range("M15").Value = -13219.8608
$ws.Range("H40").Value = 2500
$ws.Range("J40").Value = 2000
$ws.Range("L40").Value = 2000
$ws.Range("N40").Value = -2350
$ws.Range("H43").Value = 1519.3334
$ws.Range("J43").Value = 1540
$ws.Range("L43").Value = 1540
$ws.Range("N43").Value = -1678
$ws.Range("H62").Value = 2999.3333
$ws.Range("I62").Value = 2998.6667
$ws.Range("K62").Value = 2998.6667
$ws.Range("M62").Value = -2374.6667
$ws.Range("H65").Value = 2999.3333
$ws.Range("I65").Value = 2998.6667
$ws.Range("K65").Value = 14993.3335
$ws.Range("M65").Value = -11873.3335
$ws.Range("H70").Value = 15293.429
$ws.Range("I70").Value = 1286.6666
$ws.Range("K70").Value = 3859.9998
$ws.Range("M70").Value = -3589.9998
$ws.Range("H73").Value = 15293.429
$ws.Range("I73").Value = 1286.6666
$ws.Range("K73").Value = 3859.9998
$ws.Range("M73").Value = -2923.9998
$ws.Range("H116").Value = 13595.909
$ws.Range("J116").Value = 5793.5713
$ws.Range("L116").Value = 5793.5713
$ws.Range("N116").Value = -12677.5713
$ws.Range("H132").Value = 1035.7021
$ws.Range("I132").Value = 885.0714
$ws.Range("K132").Value = 2655.2142
$ws.Range("M132").Value = -125.2142000000003
$ws.Range("H135").Value = 494.2381
$ws.Range("I135").Value = 478.95
$ws.Range("K135").Value = 4310.55
$ws.Range("M135").Value = -1775.55
$ws.Range("H137").Value = 35944.863
$ws.Range("I137").Value = 1390.9048
$ws.Range("K137").Value = 4172.7144
$ws.Range("M137").Value = -1622.7144
$ws.Range("H138").Value = 3247.147
$ws.Range("I138").Value = 2965.36
$ws.Range("K138").Value = 8896.08
$ws.Range("M138").Value = -3756.08
$ws.Range("H141").Value = 935906.8
$ws.Range("I141").Value = 1168439.4
$ws.Range("K141").Value = 3505318.2
$ws.Range("M141").Value = -3500138.2

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2872.0723
$ws.Range("I32").Value = 2218.0435
$ws.Range("K32").Value = 2218.0435
$ws.Range("M32").Value = -1931.0435
$ws.Range("H45").Value = 1487
$ws.Range("I45").Value = 1013.625
$ws.Range("K45").Value = 1013.625
$ws.Range("M45").Value = -636.625
$ws.Range("H97").Value = 731.9091
$ws.Range("I97").Value = 543.875
$ws.Range("K97").Value = 543.875
$ws.Range("M97").Value = -47.875
$ws.Range("H102").Value = 1893.7368
$ws.Range("I102").Value = 1598.8
$ws.Range("K102").Value = 1598.8
$ws.Range("M102").Value = 23.20000000000005
$ws.Range("H132").Value = 1365.4615
$ws.Range("I132").Value = 909.10345
$ws.Range("J132").Value = 2688.9
$ws.Range("K132").Value = 2727.31035
$ws.Range("L132").Value = 8066.700000000001
$ws.Range("M132").Value = -197.3103499999997
$ws.Range("N132").Value = -13126.7

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 852.25
$ws.Range("I16").Value = 829.2857
$ws.Range("K16").Value = 829.2857
$ws.Range("M16").Value = -542.2857
$ws.Range("H22").Value = 956
$ws.Range("J22").Value = 1141.3334
$ws.Range("L22").Value = 1141.3334
$ws.Range("N22").Value = -1841.3334
$ws.Range("H31").Value = 2224.4167
$ws.Range("I31").Value = 1766.1666
$ws.Range("J31").Value = 2682.6667
$ws.Range("K31").Value = 1766.1666
$ws.Range("L31").Value = 2682.6667
$ws.Range("M31").Value = -1471.1666
$ws.Range("N31").Value = -3272.6667
$ws.Range("H33").Value = 2355
$ws.Range("I33").Value = 2355
$ws.Range("K33").Value = 2355
$ws.Range("M33").Value = -1976
$ws.Range("H34").Value = 2224.4167
$ws.Range("I34").Value = 1766.1666
$ws.Range("J34").Value = 2682.6667
$ws.Range("K34").Value = 1766.1666
$ws.Range("L34").Value = 2682.6667
$ws.Range("M34").Value = -1564.1666
$ws.Range("N34").Value = -3086.6667
$ws.Range("H113").Value = 852.25
$ws.Range("I113").Value = 829.2857
$ws.Range("K113").Value = 829.2857
$ws.Range("M113").Value = 1340.7143
$ws.Range("H134").Value = 1628.125
$ws.Range("I134").Value = 1363.3667
$ws.Range("K134").Value = 4090.1001
$ws.Range("M134").Value = -1555.1001

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 790.625
$ws.Range("I5").Value = 766.6667
$ws.Range("K5").Value = 2300.0001
$ws.Range("M5").Value = -2188.0001
$ws.Range("H26").Value = 569
$ws.Range("I26").Value = 569.3333
$ws.Range("K26").Value = 1707.9999
$ws.Range("M26").Value = -1419.9999
$ws.Range("H33").Value = 117.545456
$ws.Range("J33").Value = 111.25
$ws.Range("L33").Value = 667.5
$ws.Range("N33").Value = -1233.5
$ws.Range("H120").Value = 9500
$ws.Range("I120").Value = 9000
$ws.Range("J120").Value = 10000
$ws.Range("K120").Value = 27000
$ws.Range("L120").Value = 30000
$ws.Range("M120").Value = -22162
$ws.Range("N120").Value = -39676
$ws.Range("H131").Value = 814.33
$ws.Range("J131").Value = 814.47473
$ws.Range("L131").Value = 2443.42419
$ws.Range("N131").Value = -12523.42419
$ws.Range("H132").Value = 1808.1666
$ws.Range("I132").Value = 1349.75
$ws.Range("K132").Value = 12147.75
$ws.Range("M132").Value = -9617.75
$ws.Range("H135").Value = 790.625
$ws.Range("I135").Value = 766.6667
$ws.Range("K135").Value = 6900.0003
$ws.Range("M135").Value = -4365.0003
$ws.Range("H140").Value = 1516.9459
$ws.Range("J140").Value = 2318.7058
$ws.Range("L140").Value = 6956.117400000001
$ws.Range("N140").Value = -17316.1174
$ws.Range("H141").Value = 3844.3076
$ws.Range("I141").Value = 3844.3076
$ws.Range("K141").Value = 11532.9228
$ws.Range("M141").Value = -6352.9228

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4000
$ws.Range("I80").Value = 4000
$ws.Range("K80").Value = 4000
$ws.Range("M80").Value = -3002
$ws.Range("H83").Value = 4000
$ws.Range("I83").Value = 4000
$ws.Range("K83").Value = 20000
$ws.Range("M83").Value = -15008

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 12293.714
$ws.Range("J43").Value = 12293.714
$ws.Range("L43").Value = 12293.714
$ws.Range("N43").Value = -12679.714
$ws.Range("H122").Value = 8375
$ws.Range("I122").Value = 1750
$ws.Range("K122").Value = 5250
$ws.Range("M122").Value = -2800
$ws.Range("H136").Value = 2132.6128
$ws.Range("I136").Value = 1309.174
$ws.Range("K136").Value = 3927.522
$ws.Range("M136").Value = -1377.522

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 53339.266
$ws.Range("I122").Value = 87373.89
$ws.Range("K122").Value = 262121.67
$ws.Range("M122").Value = -259671.67
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("H136").Value = 19843906
$ws.Range("I136").Value = 25255042
$ws.Range("J136").Value = 3075
$ws.Range("K136").Value = 75765126
$ws.Range("L136").Value = 9225
$ws.Range("M136").Value = -75762576
$ws.Range("N136").Value = -14325
